$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window view size change (workbook.xml bookViews) ---
$excel.Width = 28800
$excel.Height = 13390

# --- Update existing strings ---
# A2 / A3 changed from "1qaz" -> "A00051", and old "www" (C-column values use "A"/"www") -> "B"
$ws.Range("A2").Value = "A00051"
$ws.Range("A3").Value = "A0005A"
$ws.Range("A4").Value = "A0006A"

# --- B column: new values, B2/B3 = "A", B4 = "B" ---
$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "A"
$ws.Range("B4").Value = "B"

# --- C column: all become "B" ---
$ws.Range("C2").Value = "B"
$ws.Range("C3").Value = "B"
$ws.Range("C4").Value = "B"

# --- D column stays "A" ---
$ws.Range("D2").Value = "A"
$ws.Range("D3").Value = "A"
$ws.Range("D4").Value = "A"

# --- G,H,I,J,K numeric values ---
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1

$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1

# --- Apply new "wrap text, centered vertically, no fill" style to column A (A2:A4) ---
$ws.Range("A2:A4").WrapText = $true
$ws.Range("A2:A4").VerticalAlignment = -4108  # xlCenter

# --- Apply green fill wrap style to G2:G4 (fillId=2, green FF00B050) ---
$ws.Range("G2:G4").Interior.Color = 0x50B000
$ws.Range("G2:G4").WrapText = $true
$ws.Range("G2:G4").VerticalAlignment = -4108

# --- selection / outline level row changes ---
$ws.Range("A5:XFD5").Select()
$ws.Outline.ShowLevels(3, 3)
